# Append 9 new "user_detail" rows (ids 110021-110029) below the existing
# data block (which currently ends at row 21), mirroring the layout/values
# of the rows already present on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(110021, 7316931025, "Magdalena Weber", "magdalena.weber@xyz.com", 932122450),
    @(110022, 9137847236, "Adrienne Hoffman", "adrienne.hoffman@xyz.com", 848488000),
    @(110023, 8428758532, "Adrienne Mcgee", "adrienne.mcgee@xyz.com", 894773246),
    @(110024, 9804209494, "Amare Coleman", "amare.coleman@xyz.com", 956554588),
    @(110025, 7105248214, "Dawson Ibarra", "dawson.ibarra@xyz.com", 765455583),
    @(110026, 9316557128, "Elvis Mcmillan", "elvis.mcmillan@xyz.com", 884282274),
    @(110027, 8103486949, "Steve George", "steve.george@xyz.com", 971073663),
    @(110028, 9601932866, "Colton Elliott", "colton.elliott@xyz.com", 809908673),
    @(110029, 9317596765, "Carolyn Rodriguez", "carolyn.rodriguez@xyz.com", 818876429)
)

$startRow = 22
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # id, uin, name, email, mobile
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]

    # status_code, lang_code, last_login_method - same constants as every
    # other row in the table
    $ws.Cells.Item($r, 6).Value = "ACT"
    $ws.Cells.Item($r, 7).Value = "eng"
    $ws.Cells.Item($r, 8).Value = "PWD"

    # is_active (boolean, left-aligned like the rest of the column)
    $ws.Cells.Item($r, 9).Value = $true
    $ws.Cells.Item($r, 9).HorizontalAlignment = -4131

    # cr_by, cr_dtimes, eff_dtimes
    $ws.Cells.Item($r, 10).Value = "superadmin"
    $ws.Cells.Item($r, 11).Value = "now()"
    $ws.Cells.Item($r, 12).Value = "now()"

    # matches the (visually invisible) "no fill" formatting carried by the
    # email column on every other data row
    $ws.Cells.Item($r, 4).Interior.Pattern = -4142
}

# Scroll the view down to the newly added rows and select the new id column
# block, matching where a user would naturally land after typing this data.
$ws.Range("A16").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A22:A30").Select()
